$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = "4.94 [1.23-13.22]"
    3  = "0.38 [0.06-3.45]"
    4  = "0.03 [0-0.81]"
    5  = "2.35 [0.34-13.73]"
    6  = "2.81 [0.4-15.01]"
    7  = "5.33 [0.49-17.58]"
    8  = "8.18 [1.99-23.7]"
    9  = "7.38 [0.55-25.52]"
    10 = "5.59 [1.5-26.81]"
    11 = "5.8 [1.3-22.53]"
    12 = "8.62 [2.06-23.86]"
    13 = "2.2 [0.63-9.63]"
    14 = "5.17 [1.24-24.56]"
    15 = "3.39 [1-10.03]"
    16 = "7.78 [1.9-17.95]"
    17 = "1.22 [0.22-2.39]"
    18 = "2.43 [0.46-10.41]"
    19 = "5.27 [0.39-19.83]"
    20 = "3.97 [0.33-15.23]"
    21 = "1.24 [0.33-6.98]"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 5).Value = $values[$row]
}
